# Apply the weekly update to the "Work Report" sheet:
#  - refresh report generation timestamp
#  - fill in Total Billed Amount / Total Line Items
#  - clear Scope ID #
#  - populate pricing for each line item (H16:H25)
#  - remove the "Point 30 / GND-MD" line item (row 26) by shifting
#    "Point 31 / CON-40-AAA-1-B" (row 27) up into it, keeping row 26's banding style
#  - remove the "Point 32 / ARM-8SF-GN-DL-C" line item entirely
#  - update the TOTAL row accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates ---------------------------------------------
$ws.Range("D5").Value2 = "Report Generated On: 08/26/2025 10:02 AM"
$ws.Range("C8").Value2 = 6609.42
$ws.Range("C9").Value2 = 11
$ws.Range("G10").Value2 = ""

# --- Line item pricing (previously all 0) ----------------------------------
$ws.Range("H16").Value2 = 350.53
$ws.Range("H17").Value2 = 676.26
$ws.Range("H18").Value2 = 198.88
$ws.Range("H19").Value2 = 1085.76
$ws.Range("H20").Value2 = 742.14
$ws.Range("H21").Value2 = 198.88
$ws.Range("H22").Value2 = 856.4400000000001
$ws.Range("H23").Value2 = 675.1799999999999
$ws.Range("H24").Value2 = 840.0599999999999
$ws.Range("H25").Value2 = 121.83

# --- Drop "Point 30 / GND-MD" (row 26): pull row 27's data up into it, ------
# keeping row 26's own (banded) formatting intact.
$ws.Range("A26:G26").Value2 = $ws.Range("A27:G27").Value2
$ws.Range("H26").Value2 = 863.46

# Now remove the now-duplicated old row 27 ("Point 31") and the old row 28
# ("Point 32 / ARM-8SF-GN-DL-C") entirely; remaining rows (incl. TOTAL) shift up.
$ws.Rows("27").Delete()
$ws.Rows("27").Delete()

# --- TOTAL row (now row 27) --------------------------------------------------
$ws.Range("H27").Value2 = 6609.420000000001
